$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Goal (per the diff): the paragraph "... Coding" (CHAPTER FOUR,
# Implementation Procedures list) turns into three numbered-list
# paragraphs:
#   1. an empty paragraph
#   2. " Coding"
#   3. " Implementation Results"   (keeps the _GoBack bookmark that used
#      to sit right after "Coding")
# ------------------------------------------------------------------

# Step 1: type the new text " Implementation Results" right before the
# hidden _GoBack bookmark that sits immediately after "Coding".
$bm = $d.Bookmarks("_GoBack")
$bmRange = $bm.Range.Duplicate
$bmRange.Collapse(1)   # wdCollapseStart
$bmRange.InsertBefore(" Implementation Results")

# Step 2: split the paragraph right after "Coding" so "Coding" and
# " Implementation Results" land in two separate paragraphs (the
# bookmark - and now " Implementation Results" - stay with the second
# half).
$afterCoding = $d.Content
$afterCoding.Find.ClearFormatting()
$afterCoding.Find.Execute("Coding") | Out-Null
$afterCoding.Collapse(0)   # wdCollapseEnd
$afterCoding.InsertParagraphAfter()

# Step 3: add a new, empty numbered-list paragraph above " Coding".
# Word's Range.InsertParagraphBefore/After treat a collapsed range that
# sits exactly on a paragraph boundary as "insert next to the whole
# paragraph" rather than as a true split at that text offset, so a
# plain SetRange(start,start)+InsertParagraphBefore on the " Coding"
# paragraph doesn't split it. Work around this by dropping a unique
# marker at the very start of that paragraph, splitting right after the
# marker (now a true interior position), and then deleting the marker.
$codingParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq " Coding`r") {
        $codingParagraph = $p
    }
}
$markStart = $codingParagraph.Range.Duplicate
$markStart.Collapse(1)   # wdCollapseStart
$markStart.InsertBefore("@@SPLITMARK@@")

$afterMark = $d.Content
$afterMark.Find.ClearFormatting()
$afterMark.Find.Execute("@@SPLITMARK@@") | Out-Null
$afterMark.Collapse(0)   # wdCollapseEnd
$afterMark.InsertParagraphAfter()

$markRange = $d.Content
$markRange.Find.ClearFormatting()
$markRange.Find.Execute("@@SPLITMARK@@") | Out-Null
$markRange.Text = ""
